$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-22"

# Update the header label for the current-year column (shared string)
$ws.Range("I1").Value = "2022 (through 06-22)"

# Update July total (row 7) and grand Total row (row 14) for the 2022 column (I)
$ws.Range("I7").Value = 104
$ws.Range("I14").Value = 767
